$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as text (inlineStr) in the source data. Several of the
# updated values look like plain numbers (e.g. "1.007"), which Excel would
# otherwise silently reinterpret as a numeric value. Pre-set those specific
# cells to Text format so the assigned string is preserved verbatim.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.540.18"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.838.22"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -2.45%  "

$ws.Range("D5").Value = "316.80"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("E6").Value = "  -2.07%  "

$ws.Range("D7").Value = "0.4304"
$ws.Range("E7").Value = "  -1.77%  "

$ws.Range("D8").Value = "0.3713"
$ws.Range("E8").Value = "  -1.86%  "

$ws.Range("D9").Value = "0.07291"
$ws.Range("E9").Value = "  -1.28%  "

$ws.Range("D10").Value = "0.8697"
$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("D11").Value = "21.23"
$ws.Range("E11").Value = "  -1.36%  "

$ws.Range("D12").Value = "1.845.04"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").Value = "6.711"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").Value = "5.377"
$ws.Range("E14").Value = "  -2.14%  "

$ws.Range("D15").Value = "0.07117"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").Value = "88.57"
$ws.Range("E16").Value = "  +4.42%  "

$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  -2.49%  "

$ws.Range("D18").Value = "0.000008952"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  -2.22%  "

$ws.Range("E20").Value = "  -0.81%  "

$ws.Range("D21").Value = "27.552.87"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").Value = "5.182"
$ws.Range("E22").Value = "  -1.76%  "

$ws.Range("D23").Value = "10.96"
$ws.Range("E23").Value = "  -2.49%  "

$ws.Range("D24").Value = "2.068.24"
$ws.Range("E24").Value = "  -1.70%  "

$ws.Range("D25").Value = "2.009"

$ws.Range("D26").Value = "154.45"
$ws.Range("E26").Value = "  -2.59%  "

$ws.Range("D27").Value = "18.53"
$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("D28").Value = "2.162"
$ws.Range("E28").Value = "  +8.82%  "

$ws.Range("D29").Value = "5.312"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "117.47"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("D31").Value = "0.08884"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "0.7704"
$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("D34").Value = "4.504"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").Value = "2.910"
$ws.Range("E35").Value = "  -3.15%  "

$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("E37").Value = "  -1.89%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01965"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.05302"
$ws.Range("E39").Value = "  +0.97%  "

$ws.Range("D40").Value = "7.168"
$ws.Range("E40").Value = "  +4.68%  "

$ws.Range("D41").Value = "2.878"
$ws.Range("E41").Value = "  +1.16%  "

$ws.Range("D42").Value = "0.5104"
$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").Value = "0.1678"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").Value = "8.730"
$ws.Range("E44").Value = "  +0.39%  "

$ws.Range("D45").Value = "10.59"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").Value = "106.75"
$ws.Range("E46").Value = "  -2.96%  "

$ws.Range("D47").Value = "0.4729"
$ws.Range("E47").Value = "  +0.90%  "

$ws.Range("D48").Value = "0.06440"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").Value = "1.008"
$ws.Range("E49").Value = "  -2.18%  "

$ws.Range("D50").Value = "1.676"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").Value = "1.839"
$ws.Range("E51").Value = "  -2.37%  "
